$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = 12
$ws.Range("E17").Value = 13
$ws.Range("E18").Value = 12
$ws.Range("E19").Value = 5
$ws.Range("E22").Value = 3
$ws.Range("E23").Value = 2

$ws.Range("I1").Select()
